$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (shifts existing A:D data to B:E)
$ws.Columns.Item(1).Insert()

# Insert a new row before row 1 (shifts existing data down by one row)
$ws.Rows.Item(1).Insert()

# New header row (row 1) for columns B:E
$ws.Cells.Item(1,2).Value = "Valid"
$ws.Cells.Item(1,3).Value = "T"
$ws.Cells.Item(1,4).Value = "Z"
$ws.Cells.Item(1,5).Value = "p-value"

# Row labels for column A (rows 2:20)
$labels = @(
  "CyclomaticComplexity(CC) & CyclomaticComplexity(CC)",
  "CyclomaticComplexity(CC) & EffortToImplement",
  "MaintainabilityIndex & MaintainabilityIndex",
  "NbUniqueOperands & NbUniqueOperands",
  "NbOperands & NbOperands",
  "NbOperands & EffortToImplement",
  "NbUniqueOperators & NbUniqueOperators",
  "NbOperators & NbOperators",
  "NbOperators & EffortToImplement",
  "ProgramLength & ProgramLength",
  "VocabularySize & VocabularySize",
  "ProgramVolume & ProgramVolume",
  "DifficultyLevel & DifficultyLevel",
  "ProgramLevel & ProgramLevel",
  "EffortToImplement & CyclomaticComplexity(CC)",
  "EffortToImplement & NbOperands",
  "EffortToImplement & NbOperators",
  "EffortToImplement & EffortToImplement",
  "TimeToImplement & TimeToImplement"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
  $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}

# Widen column A to fit the new row labels
$ws.Columns.Item(1).ColumnWidth = 53.666666666666664

Write-Host "Done"
